$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 24.7167224679393
$ws.Range("B3").Value = 24.59982259580825
$ws.Range("B4").Value = 24.53778568476006
$ws.Range("B5").Value = 24.51497430297423
$ws.Range("B6").Value = 24.51133611891744
$ws.Range("B7").Value = 24.53746802174898
$ws.Range("B8").Value = 24.67440712053321
$ws.Range("B9").Value = 25.0191449088072
$ws.Range("B10").Value = 25.31714446138998
$ws.Range("B11").Value = 25.46197415585464
$ws.Range("B12").Value = 25.51810738865101
$ws.Range("B13").Value = 25.50596146923756
$ws.Range("B14").Value = 25.46656668871028
$ws.Range("B15").Value = 25.44260282438995
$ws.Range("B16").Value = 25.30786235983332
$ws.Range("B17").Value = 25.22754735345756
$ws.Range("B18").Value = 25.18222745328625
$ws.Range("B19").Value = 25.16703449694991
$ws.Range("B20").Value = 25.23600675134699
$ws.Range("B21").Value = 25.47810325779687
$ws.Range("B22").Value = 25.64382064959544
$ws.Range("B23").Value = 25.55470358951705
$ws.Range("B24").Value = 25.23217959345524
$ws.Range("B25").Value = 24.91788983770538

$ws.Range("C2").Value = 10.13485524762043
$ws.Range("C3").Value = 9.692736206038283
$ws.Range("C4").Value = 9.414989312602454
$ws.Range("C5").Value = 9.300425307040966
$ws.Range("C6").Value = 9.281325014603052
$ws.Range("C7").Value = 9.413449567500896
$ws.Range("C8").Value = 9.983829781387236
$ws.Range("C9").Value = 11.04485168006714
$ws.Range("C10").Value = 11.78064162991518
$ws.Range("C11").Value = 12.10443109871545
$ws.Range("C12").Value = 12.22537352718984
$ws.Range("C13").Value = 12.19940207905855
$ws.Range("C14").Value = 12.11441507812274
$ws.Range("C15").Value = 12.0621379837067
$ws.Range("C16").Value = 11.75925260821422
$ws.Range("C17").Value = 11.570568998335
$ws.Range("C18").Value = 11.46101971588521
$ws.Range("C19").Value = 11.42375568738313
$ws.Range("C20").Value = 11.59076137752967
$ws.Range("C21").Value = 12.1394238447465
$ws.Range("C22").Value = 12.48823918783705
$ws.Range("C23").Value = 12.30299235583367
$ws.Range("C24").Value = 11.58163573604691
$ws.Range("C25").Value = 10.76494894192883

$ws.Range("E2").Value = 9.561315595542441
$ws.Range("E3").Value = 9.414366341636255
$ws.Range("E4").Value = 9.322784199976056
$ws.Range("E5").Value = 9.285146352752079
$ws.Range("E6").Value = 9.278878135296926
$ws.Range("E7").Value = 9.32227785686103
$ws.Range("E8").Value = 9.510942486233173
$ws.Range("E9").Value = 9.869169514381964
$ws.Range("E10").Value = 10.12391431551542
$ws.Range("E11").Value = 10.23771238428187
$ws.Range("E12").Value = 10.28048402039985
$ws.Range("E13").Value = 10.27128699699678
$ws.Range("E14").Value = 10.24123778826811
$ws.Range("E15").Value = 10.22278933507538
$ws.Range("E16").Value = 10.11643352547826
$ws.Range("E17").Value = 10.05063789777369
$ws.Range("E18").Value = 10.01259880059664
$ws.Range("E19").Value = 9.999686551438419
$ws.Range("E20").Value = 10.05766231258268
$ws.Range("E21").Value = 10.25007285129884
$ws.Range("E22").Value = 10.37394277321505
$ws.Range("E23").Value = 10.30800997320665
$ws.Range("E24").Value = 10.05448723434332
$ws.Range("E25").Value = 9.773652116430362

$ws.Range("F2").Value = 61.85349413635895
$ws.Range("F3").Value = 60.95161973743883
$ws.Range("F4").Value = 60.39682743005936
$ws.Range("F5").Value = 60.170664043623
$ws.Range("F6").Value = 60.1331101446881
$ws.Range("F7").Value = 60.39377740204621
$ws.Range("F8").Value = 61.54285112532985
$ws.Range("F9").Value = 63.78009015468538
$ws.Range("F10").Value = 65.40361354098171
$ws.Range("F11").Value = 66.13568065632846
$ws.Range("F12").Value = 66.41179021458967
$ws.Range("F13").Value = 66.3523771204014
$ws.Range("F14").Value = 66.15841955821044
$ws.Range("F15").Value = 66.0394654757996
$ws.Range("F16").Value = 65.35562560179031
$ws.Range("F17").Value = 64.93432461679328
$ws.Range("F18").Value = 64.69140354281073
$ws.Range("F19").Value = 64.60905707400441
$ws.Range("F20").Value = 64.97923620756941
$ws.Range("F21").Value = 66.21542105256043
$ws.Range("F22").Value = 67.01680138655274
$ws.Range("F23").Value = 66.5897445350214
$ws.Range("F24").Value = 64.95893387395738
$ws.Range("F25").Value = 63.17766129709469

$ws.Range("G2").Value = 3.792688912083429
$ws.Range("G3").Value = 3.797775535139773
$ws.Range("G4").Value = 3.801055111744827
$ws.Range("G5").Value = 3.80243106535808
$ws.Range("G6").Value = 3.802661932326575
$ws.Range("G7").Value = 3.801073508165697
$ws.Range("G8").Value = 3.794410433407048
$ws.Range("G9").Value = 3.782576701470478
$ws.Range("G10").Value = 3.774622370291258
$ws.Range("G11").Value = 3.771161935831185
$ws.Range("G12").Value = 3.76987409160364
$ws.Range("G13").Value = 3.77015045193393
$ws.Range("G14").Value = 3.771055533249272
$ws.Range("G15").Value = 3.771612852729306
$ws.Range("G16").Value = 3.774851682641195
$ws.Range("G17").Value = 3.776878949161209
$ws.Range("G18").Value = 3.778059864464355
$ws.Range("G19").Value = 3.778462264044786
$ws.Range("G20").Value = 3.776661603900151
$ws.Range("G21").Value = 3.770789078417329
$ws.Range("G22").Value = 3.76708238335764
$ws.Range("G23").Value = 3.769048757349729
$ws.Range("G24").Value = 3.776759817671207
$ws.Range("G25").Value = 3.785647273890865

$ws.Range("J2").Value = 12.72773707422013
$ws.Range("J3").Value = 12.62877527113188
$ws.Range("J4").Value = 12.56733002910251
$ws.Range("J5").Value = 12.54212457128045
$ws.Range("J6").Value = 12.53792935612864
$ws.Range("J7").Value = 12.56699076449482
$ws.Range("J8").Value = 12.69376000852095
$ws.Range("J9").Value = 12.93669711247255
$ws.Range("J10").Value = 13.11141074842286
$ws.Range("J11").Value = 13.19000347088513
$ws.Range("J12").Value = 13.21963097842466
$ws.Range("J13").Value = 13.21325622169893
$ws.Range("J14").Value = 13.19244366375428
$ws.Range("J15").Value = 13.17967774579651
$ws.Range("J16").Value = 13.10625615463603
$ws.Range("J17").Value = 13.06098297291457
$ws.Range("J18").Value = 13.0348600916099
$ws.Range("J19").Value = 13.02600128599844
$ws.Range("J20").Value = 13.06581101214485
$ws.Range("J21").Value = 13.19856050128486
$ws.Range("J22").Value = 13.28453675468677
$ws.Range("J23").Value = 13.2387234320264
$ws.Range("J24").Value = 13.06362855016396
$ws.Range("J25").Value = 12.87161485726535

$ws.Range("L2").Value = 9.9377712289069
$ws.Range("L3").Value = 9.995897139369156
$ws.Range("L4").Value = 10.03350031881593
$ws.Range("L5").Value = 10.04930677760731
$ws.Range("L6").Value = 10.05196064032718
$ws.Range("L7").Value = 10.03371153294739
$ws.Range("L8").Value = 9.95741686557553
$ws.Range("L9").Value = 9.822911945388444
$ws.Range("L10").Value = 9.733197465744684
$ws.Range("L11").Value = 9.694339327514742
$ws.Range("L12").Value = 9.679903993112484
$ws.Range("L13").Value = 9.683000497978002
$ws.Range("L14").Value = 9.693146132820797
$ws.Range("L15").Value = 9.6993969689127
$ws.Range("L16").Value = 9.735776115627687
$ws.Range("L17").Value = 9.758592788949686
$ws.Range("L18").Value = 9.771900290475131
$ws.Range("L19").Value = 9.77643762224066
$ws.Range("L20").Value = 9.756144887478314
$ws.Range("L21").Value = 9.690158542980944
$ws.Range("L22").Value = 9.648660522480917
$ws.Range("L23").Value = 9.67066033751415
$ws.Range("L24").Value = 9.757250992148322
$ws.Range("L25").Value = 9.857692458735089

$ws.Range("M2").Value = 19.10156422396014
$ws.Range("M3").Value = 19.18990972600128
$ws.Range("M4").Value = 19.24975020337086
$ws.Range("M5").Value = 19.27553947231622
$ws.Range("M6").Value = 19.27990644674785
$ws.Range("M7").Value = 19.25009232701553
$ws.Range("M8").Value = 19.13086229630577
$ws.Range("M9").Value = 18.94160389344235
$ws.Range("M10").Value = 18.82991885427995
$ws.Range("M11").Value = 18.7850946417649
$ws.Range("M12").Value = 18.76898447374727
$ws.Range("M13").Value = 18.77241561295433
$ws.Range("M14").Value = 18.78375191662102
$ws.Range("M15").Value = 18.79080832612991
$ws.Range("M16").Value = 18.8329689024102
$ws.Range("M17").Value = 18.8603677407821
$ws.Range("M18").Value = 18.87668965989719
$ws.Range("M19").Value = 18.88231253692908
$ws.Range("M20").Value = 18.85739281035193
$ws.Range("M21").Value = 18.7803986981712
$ws.Range("M22").Value = 18.7351150242127
$ws.Range("M23").Value = 18.7588217717626
$ws.Range("M24").Value = 18.85873600136892
$ws.Range("M25").Value = 18.98801253705334
